$wb = $excel.ActiveWorkbook

$wsFull = $wb.Worksheets.Item("Full results")
$wsPlot = $wb.Worksheets.Item("For plotting")

# --- Sheet "Full results" ---
# Row 2 (health_pc / NULL MODEL)
$wsFull.Range("C2").Value = 0.908941580497606
$wsFull.Range("D2").Value = 0.0911546741138437
$wsFull.Range("E2").Value = 1.00009625461145
$wsFull.Range("J2").Value = 0.0911459009005673
$wsFull.Range("K2").Value = 0.09185541392866
$wsFull.Range("L2").Value = -0.0118209738535919
$wsFull.Range("M2").Value = 0.0062307079296351
$wsFull.Range("N2").Value = 0.0800344400750681

# Row 3 (health_pc / CONDITIONAL MODEL)
$wsFull.Range("F3").Value = 0.890888161156763
$wsFull.Range("G3").Value = 0.0918642554358372

# Row 4 (health_pc / COMPLETE MODEL)
$wsFull.Range("H4").Value = 0.9027102728336
$wsFull.Range("I4").Value = 0.0768200707292788
$wsFull.Range("O4").Value = 0.0973766088302024

# --- Sheet "For plotting" ---
# Row 2 (w / health_pc)
$wsPlot.Range("C2").Value = 0.0911459009005673
$wsPlot.Range("D2").Value = 0.0360600100989626
$wsPlot.Range("E2").Value = 0.146231791702172
$wsPlot.Range("F2").Value = 948

# Row 3 (IOLIB / health_pc)
$wsPlot.Range("C3").Value = 0.0800344400750681
$wsPlot.Range("D3").Value = 0.0207353841066964
$wsPlot.Range("E3").Value = 0.13933349604344
$wsPlot.Range("F3").Value = 948

# Row 4 (IORAD / health_pc)
$wsPlot.Range("C4").Value = 0.0973766088302024
$wsPlot.Range("D4").Value = 0.01605405064388
$wsPlot.Range("E4").Value = 0.178699167016525
$wsPlot.Range("F4").Value = 948
